$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF) - copy format from H1 so the
# same shared cell style (bold, centered, bordered) is reused.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @{
  2  = @(6, 6)
  3  = @(8, 8)
  4  = @(8, 9)
  5  = @(7, 7)
  6  = @(9, 9)
  7  = @(9, 9)
  8  = @(6, 6)
  9  = @(6, 7)
  10 = @(9, 9)
  11 = @(1, 1)
  12 = @(8, 8)
  13 = @(7, 7)
  14 = @(8, 8)
  15 = @(8, 8)
  16 = @(9, 9)
  17 = @(8, 8)
  18 = @(10, 10)
  19 = @(8, 8)
  20 = @(7, 7)
  21 = @(6, 7)
  22 = @(8, 8)
  23 = @(7, 7)
  24 = @(7, 7)
  25 = @(7, 7)
  26 = @(8, 8)
  27 = @(8, 8)
  28 = @(8, 8)
  29 = @(7, 7)
  30 = @(7, 7)
  31 = @(7, 7)
  32 = @(7, 7)
  33 = @(7, 7)
  34 = @(7, 7)
  35 = @(6, 7)
  36 = @(6, 6)
  37 = @(8, 8)
  38 = @(8, 8)
  39 = @(8, 8)
  40 = @(7, 7)
  41 = @(7, 7)
  42 = @(8, 8)
  43 = @(8, 8)
  44 = @(6, 7)
  45 = @(9, 9)
  46 = @(5, 6)
  47 = @(6, 6)
  48 = @(8, 8)
  49 = @(7, 7)
  50 = @(8, 8)
  51 = @(7, 7)
  52 = @(1, 2)
  53 = @(7, 7)
  54 = @(9, 9)
  55 = @(8, 8)
  56 = @(8, 8)
  57 = @(7, 8)
  58 = @(9, 9)
  59 = @(8, 8)
  60 = @(6, 7)
  61 = @(5, 5)
  62 = @(3, 3)
  63 = @(5, 5)
}

foreach ($r in ($data.Keys | Sort-Object)) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
